# Actualización automática del mapa: agrega el nuevo caso (fila 94) a la
# hoja "AYKO" con los datos del reclamo 6497.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

# Algunas columnas contienen valores que parecen numéricos/fechas
# (p.ej. "6497", "7/25/2025", "14", "808533127") pero deben preservarse
# como texto, igual que el resto de la hoja. Se fuerza el formato de
# texto antes de escribir y se limpia el formato después para no dejar
# estilos adicionales en las celdas.
$ws.Range("A$row`:E$row").NumberFormat = "@"

$ws.Cells.Item($row, 1).Value  = "6497"
$ws.Cells.Item($row, 2).Value  = "7/25/2025"
$ws.Cells.Item($row, 3).Value  = "SEGUI, JUAN FRANCISCO 4691"
$ws.Cells.Item($row, 4).Value  = "14"
$ws.Cells.Item($row, 5).Value  = "808533127"
$ws.Cells.Item($row, 6).Value  = "AYKO"
$ws.Cells.Item($row, 7).Value  = "Pendiente"
$ws.Cells.Item($row, 8).Value  = "Inclinada"
$ws.Cells.Item($row, 9).Value  = 1
$ws.Cells.Item($row, 10).Value = "Aplomo"
$ws.Cells.Item($row, 11).Value = "Sin equipos"
$ws.Cells.Item($row, 12).Value = "Pasante"
$ws.Cells.Item($row, 13).Value = -58.422229
$ws.Cells.Item($row, 14).Value = -34.573148
$ws.Cells.Item($row, 15).Value = "Palermo"
$ws.Cells.Item($row, 16).Value = "Capital Sur"

$ws.Range("A$row`:P$row").ClearFormats()
